# Generate Report for Handback
# Update the timestamp cells (stored as text) in the three sheets to reflect
# the newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G)
$overview.Range("G2").Value = "2016-09-05 23:19:30"

# zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K)
$zhcn.Range("H2").Value = "2016-09-05 23:19:25"
$zhcn.Range("K2").Value = "2016-09-05 23:19:43"

# de-de sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K)
$dede.Range("H2").Value = "2016-09-05 23:19:30"
$dede.Range("K2").Value = "2016-09-05 23:19:51"
